$d = $word.ActiveDocument

$pairs = @(
    @("343×6=2058", "376×7=2632"),
    @("760×8=6080", "787×8=6296"),
    @("633×5=3165", "836×3=2508"),
    @("938×3=2814", "573×3=1719"),
    @("367×4=1468", "106×8=848"),
    @("224×9=2016", "326×6=1956"),
    @("475×7=3325", "592×9=5328"),
    @("664×5=3320", "757×4=3028"),
    @("267×3=801",  "609×5=3045"),
    @("864×4=3456", "901×7=6307"),
    @("441×7=3087", "188×3=564"),
    @("519×5=2595", "740×4=2960"),
    @("460×8=3680", "587×7=4109"),
    @("826×7=5782", "579×9=5211"),
    @("668×4=2672", "475×5=2375"),
    @("128×8=1024", "159×5=795"),
    @("848×2=1696", "376×6=2256"),
    @("332×5=1660", "743×4=2972"),
    @("255×4=1020", "847×3=2541"),
    @("530×7=3710", "881×2=1762"),
    @("879×5=4395", "224×7=1568"),
    @("597×8=4776", "698×8=5584"),
    @("499×7=3493", "161×5=805"),
    @("429×9=3861", "580×7=4060"),
    @("831×7=5817", "151×4=604")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
